# The original document has a single hidden "_GoBack" bookmark (Word's
# automatic "last edit location" marker) sitting just after the "V" in
# "(Version 1.0)". The edit being reproduced:
#   1. Changes the quoted folder name "testing data" -> "test examples"
#      inside the sentence "...available in the folder- "testing data".".
#   2. As a natural consequence of that being the last text edit made in
#      the document, Word relocates the hidden _GoBack bookmark to sit
#      right after the newly typed text (and before the trailing period),
#      removing it from its old location near "(Version 1.0)".

$d = $word.ActiveDocument

# Locate the quoted phrase (including the smart/curly quotes used in the
# document) that needs to be retyped.
$target = $d.Content
$found = $target.Find.Execute(
    [char]0x201C + "testing data" + [char]0x201D,
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the folder name to replace"
}

# Pin the boundary with the preceding run ("Please note that t") so that
# the upcoming text replacement does not get coalesced backwards into it.
$beforeAnchor = $d.Range($target.Start, $target.Start)
$d.Bookmarks.Add("ZZTMP_BOUNDARY", $beforeAnchor)

# Re-find the phrase (collection/range handles can go stale after the
# structural edit above) and retype it, exactly like a user selecting
# the quoted text and typing the replacement over it.
$target = $d.Content
$target.Find.Execute(
    [char]0x201C + "testing data" + [char]0x201D,
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$target.Text = [char]0x201C + "test examples" + [char]0x201D

# Word drops its hidden _GoBack bookmark at the end of the text that was
# just typed. Adding a bookmark with that name automatically replaces any
# existing one elsewhere in the document (matching real Word behaviour).
$goBackPoint = $d.Range($target.End, $target.End)
$d.Bookmarks.Add("_GoBack", $goBackPoint)

# Clean up the temporary boundary marker now that the edit is complete.
$d.Bookmarks.Item("ZZTMP_BOUNDARY").Delete()
